$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.944.30'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.336.58'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.02'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.30'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.502'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.08'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.48%  '
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.70'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -4.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.121'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.58%  '
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.361.75'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.876.60'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.05'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.18'
$ws.Range("D20").ClearFormats()
$ws.Range("E21").Value = '  -1.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.82'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.50'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.21'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.88%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.41'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.50'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.11'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.38'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.70%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0738'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +4.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.22'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.38'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("E36").Value = '  -1.12%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '126.09'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -23.35%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.82'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.90%  '
$ws.Range("E39").Value = '  -0.57%  '
$ws.Range("E40").Value = '  -1.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.10'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +14.35%  '
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.941.79'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0282'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.16'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.56%  '
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("E47").Value = '  -3.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.568.65'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.77'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.51'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.69%  '
